$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2450.6
$ws.Range("J17").Value = 2450.6
$ws.Range("L17").Value = 7351.799999999999
$ws.Range("N17").Value = -7687.799999999999
$ws.Range("H55").Value = 875
$ws.Range("I55").Value = 750
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 750
$ws.Range("L55").Value = 1000
$ws.Range("M55").Value = -536
$ws.Range("N55").Value = -1428
$ws.Range("H59").Value = 2968.8
$ws.Range("I59").Value = 1955.6666
$ws.Range("K59").Value = 5866.9998
$ws.Range("M59").Value = -5309.9998
$ws.Range("H94").Value = 2018.7
$ws.Range("J94").Value = 3400
$ws.Range("L94").Value = 3400
$ws.Range("N94").Value = -4302
$ws.Range("H98").Value = 1622.475
$ws.Range("I98").Value = 1549.7646
$ws.Range("K98").Value = 1549.7646
$ws.Range("M98").Value = -51.76459999999997
$ws.Range("H100").Value = 6221.3335
$ws.Range("I100").Value = 16666
$ws.Range("J100").Value = 999
$ws.Range("K100").Value = 16666
$ws.Range("L100").Value = 999
$ws.Range("M100").Value = -16125
$ws.Range("N100").Value = -2081
$ws.Range("H101").Value = 1864.3334
$ws.Range("I101").Value = 2262.8333
$ws.Range("J101").Value = 1465.8334
$ws.Range("K101").Value = 6788.499899999999
$ws.Range("L101").Value = 4397.5002
$ws.Range("M101").Value = -5166.499899999999
$ws.Range("N101").Value = -7641.5002
$ws.Range("H113").Value = 4292.04
$ws.Range("I113").Value = 3440.15
$ws.Range("J113").Value = 7699.6
$ws.Range("K113").Value = 3440.15
$ws.Range("L113").Value = 7699.6
$ws.Range("M113").Value = -186.1500000000001
$ws.Range("N113").Value = -14207.6
$ws.Range("H122").Value = 1622.475
$ws.Range("I122").Value = 1549.7646
$ws.Range("K122").Value = 4649.293799999999
$ws.Range("M122").Value = -2199.293799999999
$ws.Range("H135").Value = 1125.4359
$ws.Range("I135").Value = 1088.4839
$ws.Range("J135").Value = 1268.625
$ws.Range("K135").Value = 9796.355099999999
$ws.Range("L135").Value = 11417.625
$ws.Range("M135").Value = -7261.355099999999
$ws.Range("N135").Value = -16487.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1136.6364
$ws.Range("I2").Value = 1282.4814
$ws.Range("K2").Value = 1282.4814
$ws.Range("M2").Value = -1169.4814
$ws.Range("H110").Value = 260087.86
$ws.Range("I110").Value = 451363.5
$ws.Range("K110").Value = 451363.5
$ws.Range("M110").Value = -449318.5
$ws.Range("H116").Value = 1136.6364
$ws.Range("I116").Value = 1282.4814
$ws.Range("K116").Value = 1282.4814
$ws.Range("M116").Value = 1011.5186
$ws.Range("H132").Value = 1422.8536
$ws.Range("I132").Value = 1541.9429
$ws.Range("K132").Value = 4625.8287
$ws.Range("M132").Value = -2095.8287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1136.6364
$ws.Range("I3").Value = 1282.4814
$ws.Range("K3").Value = 1282.4814
$ws.Range("M3").Value = -1168.4814
$ws.Range("H94").Value = 76349.414
$ws.Range("I94").Value = 956.7143
$ws.Range("K94").Value = 956.7143
$ws.Range("M94").Value = -505.7143
$ws.Range("H99").Value = 1517.75
$ws.Range("I99").Value = 1057.8667
$ws.Range("J99").Value = 2897.4
$ws.Range("K99").Value = 1057.8667
$ws.Range("L99").Value = 2897.4
$ws.Range("M99").Value = 440.1333
$ws.Range("N99").Value = -5893.4
$ws.Range("N112").ClearContents()
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("H115").Value = 33750.25
$ws.Range("J115").Value = 33750.25
$ws.Range("L115").Value = 33750.25
$ws.Range("N115").Value = -36884.25
$ws.Range("H120").Value = 59998
$ws.Range("J120").Value = 59998
$ws.Range("L120").Value = 59998
$ws.Range("N120").Value = -69674
$ws.Range("H134").Value = 1292.5227
$ws.Range("I134").Value = 1347.4054
$ws.Range("J134").Value = 1002.4286
$ws.Range("K134").Value = 4042.2162
$ws.Range("L134").Value = 3007.2858
$ws.Range("M134").Value = -1507.2162
$ws.Range("N134").Value = -8077.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1807.7715
$ws.Range("I31").Value = 1493.5416
$ws.Range("K31").Value = 1493.5416
$ws.Range("M31").Value = -1198.5416
$ws.Range("H34").Value = 1807.7715
$ws.Range("I34").Value = 1493.5416
$ws.Range("K34").Value = 1493.5416
$ws.Range("M34").Value = -1291.5416
$ws.Range("H50").Value = 38076.867
$ws.Range("J50").Value = 40088.69
$ws.Range("L50").Value = 40088.69
$ws.Range("N50").Value = -41338.69
$ws.Range("H141").Value = 169050.53
$ws.Range("J141").Value = 169050.53
$ws.Range("L141").Value = 169050.53
$ws.Range("N141").Value = -179410.53

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 1992
$ws.Range("J46").Value = 1992
$ws.Range("L46").Value = 5976
$ws.Range("N46").Value = -6158
$ws.Range("H131").Value = 3532.3076
$ws.Range("I131").Value = 2329.6667
$ws.Range("J131").Value = 4563.143
$ws.Range("K131").Value = 6989.000100000001
$ws.Range("L131").Value = 13689.429
$ws.Range("M131").Value = -1949.000100000001
$ws.Range("N131").Value = -23769.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 6638.625
$ws.Range("J13").Value = 5516.5
$ws.Range("L13").Value = 5516.5
$ws.Range("N13").Value = -5794.5
$ws.Range("H97").Value = 40000492
$ws.Range("I97").Value = 47619416
$ws.Range("K97").Value = 47619416
$ws.Range("M97").Value = -47618920
$ws.Range("H102").Value = 19286.645
$ws.Range("I102").Value = 24567.93
$ws.Range("K102").Value = 24567.93
$ws.Range("M102").Value = -22945.93
$ws.Range("H122").Value = 67451.78999999999
$ws.Range("J122").Value = 3682.6
$ws.Range("L122").Value = 11047.8
$ws.Range("N122").Value = -15947.8
$ws.Range("H126").Value = 46435.57
$ws.Range("J126").Value = 3773
$ws.Range("L126").Value = 11319
$ws.Range("N126").Value = -16259
$ws.Range("H132").Value = 2501.1365
$ws.Range("I132").Value = 2647.1538
$ws.Range("K132").Value = 7941.4614
$ws.Range("M132").Value = -5411.4614
$ws.Range("H139").Value = 78240.75
$ws.Range("J139").Value = 78240.75
$ws.Range("L139").Value = 78240.75
$ws.Range("N139").Value = -88520.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 38444
$ws.Range("J92").Value = 38444
$ws.Range("L92").Value = 38444
$ws.Range("N92").Value = -43436
$ws.Range("H132").Value = 3397.0715
$ws.Range("I132").Value = 2968.95
$ws.Range("J132").Value = 4467.375
$ws.Range("K132").Value = 8906.849999999999
$ws.Range("L132").Value = 13402.125
$ws.Range("M132").Value = -6376.849999999999
$ws.Range("N132").Value = -18462.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M25").ClearContents()
$ws.Range("H25").Value = 100026
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 100026
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 100026
$ws.Range("N25").Value = -100612
$ws.Range("H28").Value = 20000
$ws.Range("J28").Value = 20000
$ws.Range("L28").Value = 20000
$ws.Range("N28").Value = -20696
$ws.Range("H58").Value = 8713.700000000001
$ws.Range("I58").Value = 8713.700000000001
$ws.Range("K58").Value = 8713.700000000001
$ws.Range("M58").Value = -8405.700000000001
$ws.Range("H81").Value = 1474.5
$ws.Range("I81").Value = 1299.3334
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 2598.6668
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -1537.6668
$ws.Range("N81").Value = -6122
$ws.Range("H84").Value = 1474.5
$ws.Range("I84").Value = 1299.3334
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 12993.334
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -7689.333999999999
$ws.Range("N84").Value = -30608
$ws.Range("H100").Value = 896.10345
$ws.Range("I100").Value = 772.36365
$ws.Range("K100").Value = 1544.7273
$ws.Range("M100").Value = -1003.7273
$ws.Range("H126").Value = 1991.65
$ws.Range("I126").Value = 1991.0555
$ws.Range("K126").Value = 5973.166499999999
$ws.Range("M126").Value = -3503.166499999999
